$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = "2022-07-21 20:57:28"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6753975"
$ws.Range("B3").Value = "Severin Standgrill mit Grillplatte PG 8563"
$ws.Range("C3").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-standgrill-mit-grillplatte-pg-8563/p/6753975"
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "Severin"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "74.50"
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = $null
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N3").Value = "Severin Standgrill mit Grillplatte PG 8563 50% Aktion 74.50 Schweizer Franken statt 149.00 Schweizer Franken"
$ws.Range("O3").Value = "2022-07-21 20:57:28"

# Row 4
$ws.Range("O4").Value = "2022-07-21 20:57:28"

# Row 5
$ws.Range("O5").Value = "2022-07-21 20:57:28"

# Row 6
$ws.Range("O6").Value = "2022-07-21 20:57:28"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "3591269"
$ws.Range("B7").Value = "Varta Longlife Max Power C 2er Bli"
$ws.Range("C7").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-c-2er-bli/p/3591269"
$ws.Range("D7").Value = "2ST"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 5
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "8.95"
$ws.Range("I7").Value = "4.48/1ST"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "4.48"
$ws.Range("N7").Value = "Varta Longlife Max Power C 2er Bli 8.95 Schweizer Franken"
$ws.Range("O7").Value = "2022-07-21 20:57:28"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "6266870"
$ws.Range("B8").Value = "Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück"
$ws.Range("C8").Value = "/de/haushalt-tier/elektroartikel-batterien/reer-steckdosen-sicherung-kappe-schweiz-8-stueck/p/6266870"
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = $null
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "reer"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "3.95"
$ws.Range("I8").Value = $null
$ws.Range("J8").Value = $null
$ws.Range("K8").Value = $null
$ws.Range("L8").Value = $null
$ws.Range("M8").Value = "['haushalt-tier', 'elektroartikel-batterien']"
$ws.Range("N8").Value = "Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück 3.95 Schweizer Franken"
$ws.Range("O8").Value = "2022-07-21 20:57:28"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "5872158"
$ws.Range("B9").Value = "satrap Tischventilator Venti 1"
$ws.Range("C9").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-1/p/5872158"
$ws.Range("G9").Value = "satrap"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "29.95"
$ws.Range("M9").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N9").Value = "satrap Tischventilator Venti 1 29.95 Schweizer Franken"
$ws.Range("O9").Value = "2022-07-21 20:57:28"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "3494230"
$ws.Range("B10").Value = "Varta Electronics V13GS / V357 1er Bli"
$ws.Range("C10").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13gs--v357-1er-bli/p/3494230"
$ws.Range("D10").Value = "1ST"
$ws.Range("G10").Value = "Varta"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "4.95"
$ws.Range("I10").Value = "4.95/1ST"
$ws.Range("J10").Value = "Preis pro 1 Stück"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "4.95"
$ws.Range("L10").Value = "1ST"
$ws.Range("M10").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("N10").Value = "Varta Electronics V13GS / V357 1er Bli 4.95 Schweizer Franken"
$ws.Range("O10").Value = "2022-07-21 20:57:28"

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "6872591"
$ws.Range("B11").Value = "Bosch Küchenmaschine MUM58243 1000W"
$ws.Range("C11").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/bosch-kuechenmaschine-mum58243-1000w/p/6872591"
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "Bosch"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "239.50"
$ws.Range("M11").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N11").Value = "Bosch Küchenmaschine MUM58243 1000W 50% Aktion 239.50 Schweizer Franken statt 479.00 Schweizer Franken"
$ws.Range("O11").Value = "2022-07-21 20:57:28"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "6735643"
$ws.Range("B12").Value = "LED 31V Anschlussset Transf.+Verl.kabel"
$ws.Range("C12").Value = "/de/haushalt-tier/haushalt-kueche/uebrige-haushaltsartikel/led-31v-anschlussset-transfverlkabel/p/6735643"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = "Coop"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "9.95"
$ws.Range("M12").Value = "['haushalt-tier', 'haushalt-kueche', 'uebrige-haushaltsartikel']"
$ws.Range("N12").Value = "LED 31V Anschlussset Transf.+Verl.kabel 50% Aktion 9.95 Schweizer Franken statt 19.95 Schweizer Franken"
$ws.Range("O12").Value = "2022-07-21 20:57:28"

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "5882124"
$ws.Range("B13").Value = "Philips Avent Audio Monitors DECT-Babyphone"
$ws.Range("C13").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/philips-avent-audio-monitors-dect-babyphone/p/5882124"
$ws.Range("G13").Value = "Avent"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "99.90"
$ws.Range("M13").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete']"
$ws.Range("N13").Value = "Philips Avent Audio Monitors DECT-Babyphone 99.90 Schweizer Franken"
$ws.Range("O13").Value = "2022-07-21 20:57:28"

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "6689619"
$ws.Range("B14").Value = "Philips Dampfstation HI5919/31"
$ws.Range("C14").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/philips-dampfstation-hi591931/p/6689619"
$ws.Range("E14").Value = $null
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "Philips"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "69.50"
$ws.Range("M14").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'haushaltsgeraete-kabel']"
$ws.Range("N14").Value = "Philips Dampfstation HI5919/31 50% Aktion 69.50 Schweizer Franken statt 139.00 Schweizer Franken"
$ws.Range("O14").Value = "2022-07-21 20:57:28"

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "6425996"
$ws.Range("B15").Value = "satrap espresso XA Kolbenkaffeemaschine"
$ws.Range("C15").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-espresso-xa-kolbenkaffeemaschine/p/6425996"
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 4
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "249.00"
$ws.Range("N15").Value = "satrap espresso XA Kolbenkaffeemaschine 249.00 Schweizer Franken"
$ws.Range("O15").Value = "2022-07-21 20:57:28"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "6125818"
$ws.Range("B16").Value = "satrap Mano XA Handmixer"
$ws.Range("C16").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mano-xa-handmixer/p/6125818"
$ws.Range("N16").Value = "satrap Mano XA Handmixer 49.95 Schweizer Franken"
$ws.Range("O16").Value = "2022-07-21 20:57:28"

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "5831402"
$ws.Range("B17").Value = "Satrap Mikrowelle Micro M2"
$ws.Range("C17").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mikrowelle-micro-m2/p/5831402"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "49.95"
$ws.Range("M17").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'kuechengeraete']"
$ws.Range("N17").Value = "Satrap Mikrowelle Micro M2 50% Aktion 49.95 Schweizer Franken statt 99.90 Schweizer Franken"
$ws.Range("O17").Value = "2022-07-21 20:57:28"

# Row 18
$ws.Range("O18").Value = "2022-07-21 20:57:28"

# Row 19
$ws.Range("O19").Value = "2022-07-21 20:57:28"

# Row 20
$ws.Range("O20").Value = "2022-07-21 20:57:28"

# Row 21
$ws.Range("O21").Value = "2022-07-21 20:57:28"

# Row 22
$ws.Range("O22").Value = "2022-07-21 20:57:28"

# Row 23
$ws.Range("O23").Value = "2022-07-21 20:57:28"

# Row 24
$ws.Range("O24").Value = "2022-07-21 20:57:28"

# Row 25
$ws.Range("O25").Value = "2022-07-21 20:57:28"

# Row 26
$ws.Range("O26").Value = "2022-07-21 20:57:28"

# Row 27
$ws.Range("O27").Value = "2022-07-21 20:57:28"

# Row 28
$ws.Range("O28").Value = "2022-07-21 20:57:28"

# Row 29
$ws.Range("O29").Value = "2022-07-21 20:57:28"

# Row 30
$ws.Range("O30").Value = "2022-07-21 20:57:28"
